$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 87

# Copy formatting from the row above so the new row matches the sheet's
# existing data-row style (border/alignment), then fill in the values.
$ws.Range("A86:H86").Copy()
$ws.Range("A87:H87").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = "2025-09-02 06:49:04 UTC"
$ws.Cells.Item($row, 2).Value = "2025-09-02 12:19:04 IST"
$ws.Cells.Item($row, 3).Value = "UPDATED"
$ws.Cells.Item($row, 4).Value = "New circular processed."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"
$ws.Cells.Item($row, 6).Value = "INGOT-01-09-2025.pdf"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = 5
